$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (Price / Volume) remain stored as text, matching the source data.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"

# Apply the updated cell values from the data refresh.
$ws.Range("D2").Value = "278.65"
$ws.Range("E2").Value = "6.70%"
$ws.Range("D3").Value = "27.30"
$ws.Range("E3").Value = "0.74%"
$ws.Range("D4").Value = "4.816"
$ws.Range("E4").Value = "2.55%"
$ws.Range("D5").Value = "0.06276"
$ws.Range("E5").Value = "0.88%"
$ws.Range("D6").Value = "6.859"
$ws.Range("E6").Value = "1.68%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "3.268"
$ws.Range("E7").Value = "2.93%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "0.8776"
$ws.Range("E8").Value = "2.99%"
$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D9").Value = "0.9497"
$ws.Range("E9").Value = "4.21%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1455"
$ws.Range("E10").Value = "4.04%"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "0.05209"
$ws.Range("E11").Value = "10.40%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.07327"
$ws.Range("E12").Value = "3.34%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03133"
$ws.Range("E13").Value = "0.10%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09049"
$ws.Range("E14").Value = "-0.03%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001555"
$ws.Range("E15").Value = "1.34%"
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D16").Value = "0.0006263"
$ws.Range("E16").Value = "1.51%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "0.006104"
$ws.Range("E17").Value = "-0.44%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "3.461"
$ws.Range("E18").Value = "0.29%"
$ws.Range("D19").Value = "2.246"
$ws.Range("E19").Value = "3.70%"
$ws.Range("E20").Value = "-0.62%"
$ws.Range("D21").Value = "0.1309"
$ws.Range("E21").Value = "-0.15%"
$ws.Range("D22").Value = "3.842"
$ws.Range("E22").Value = "-5.86%"
$ws.Range("D23").Value = "0.04323"
$ws.Range("E23").Value = "1.71%"
$ws.Range("D24").Value = "0.001173"
$ws.Range("E24").Value = "-3.45%"
$ws.Range("D25").Value = "0.004279"
$ws.Range("E25").Value = "4.60%"
$ws.Range("D26").Value = "0.0001199"
$ws.Range("E26").Value = "-0.18%"
$ws.Range("D27").Value = "0.0001684"
$ws.Range("E27").Value = "2.72%"
$ws.Range("D40").Value = "0.04029"
$ws.Range("E40").Value = "3.22%"
$ws.Range("D41").Value = "0.006696"
$ws.Range("E41").Value = "62.31%"
$ws.Range("D42").Value = "0.1154"
$ws.Range("E42").Value = "3.78%"
$ws.Range("D43").Value = "0.002099"
$ws.Range("E43").Value = "-5.08%"
$ws.Range("D44").Value = "0.01380"
$ws.Range("E44").Value = "-0.66%"
$ws.Range("D45").Value = "0.00005167"
$ws.Range("E45").Value = "1.09%"
$ws.Range("E46").Value = "-0.21%"
$ws.Range("D47").Value = "2.316"
$ws.Range("E47").Value = "1,388.21%"
$ws.Range("E48").Value = "-12.17%"
$ws.Range("D49").Value = "0.00002096"
$ws.Range("E49").Value = "-0.21%"
$ws.Range("D50").Value = "0.0001996"
$ws.Range("E50").Value = "-0.21%"
